$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $NewValue
    $r.Style = "Normal"
}

$sub3 = [string][char]0x2083

Set-TextValue "D2" ('51.773.11')
Set-TextValue "E2" ('  +0.41%  ')
Set-TextValue "D3" ('2.983.74')
Set-TextValue "E3" ('  +1.52%  ')
Set-TextValue "E4" ('  -0.02%  ')
Set-TextValue "D5" ('380.52')
Set-TextValue "E5" ('  +6.41%  ')
Set-TextValue "D6" ('105.33')
Set-TextValue "E6" ('  +0.37%  ')
Set-TextValue "D7" ('0.546')
Set-TextValue "E7" ('  -0.11%  ')
Set-TextValue "E8" ('  -0.02%  ')
Set-TextValue "D9" ('0.601')
Set-TextValue "E9" ('  +0.89%  ')
Set-TextValue "D10" ('37.68')
Set-TextValue "E10" ('  +0.81%  ')
Set-TextValue "E11" ('  -0.38%  ')
Set-TextValue "D12" ('0.0848')
Set-TextValue "E12" ('  +0.55%  ')
Set-TextValue "D14" ('3.458.96')
Set-TextValue "E14" ('  +1.72%  ')
Set-TextValue "D15" ('7.53')
Set-TextValue "E15" ('  +1.39%  ')
Set-TextValue "D16" ('2.978.27')
Set-TextValue "E16" ('  +1.46%  ')
Set-TextValue "D17" ('0.963')
Set-TextValue "E17" ('  -1.62%  ')
Set-TextValue "D18" ('51.864.58')
Set-TextValue "E19" ('  +5.54%  ')
Set-TextValue "D20" ('7.44')
Set-TextValue "E20" ('  +2.10%  ')
Set-TextValue "D21" ('13.19')
Set-TextValue "E21" ('  +0.30%  ')
Set-TextValue "D22" ('0.0' + $sub3 + '0962')
Set-TextValue "E22" ('  +0.85%  ')
Set-TextValue "E23" ('  -0.01%  ')
Set-TextValue "D24" ('263.95')
Set-TextValue "E24" ('  +0.10%  ')
Set-TextValue "D25" ('2.81')
Set-TextValue "E25" ('  +4.41%  ')
Set-TextValue "D26" ('7.40')
Set-TextValue "E26" ('  +18.93%  ')
Set-TextValue "D27" ('0.172')
Set-TextValue "E27" ('  -1.41%  ')
Set-TextValue "E28" ('  -3.50%  ')
Set-TextValue "D29" ('7.51')
Set-TextValue "E29" ('  +4.90%  ')
Set-TextValue "B30" ('Dai')
Set-TextValue "C30" ('https://coinranking.com/coin/MoTuySvg7+dai-dai')
Set-TextValue "D30" ('1.00')
Set-TextValue "E30" ('  +0.01%  ')
Set-TextValue "B31" ('EthereumClassic')
Set-TextValue "C31" ('https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc')
Set-TextValue "D31" ('26.16')
Set-TextValue "E31" ('  -1.02%  ')
Set-TextValue "E32" ('  -1.98%  ')
Set-TextValue "D33" ('9.97')
Set-TextValue "E33" ('  -0.56%  ')
Set-TextValue "B34" ('OKB')
Set-TextValue "C34" ('https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb')
Set-TextValue "D34" ('51.71')
Set-TextValue "E34" ('  +1.80%  ')
Set-TextValue "B35" ('InjectiveProtocol')
Set-TextValue "C35" ('https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj')
Set-TextValue "D35" ('34.68')
Set-TextValue "E35" ('  -1.90%  ')
Set-TextValue "E36" ('  -3.89%  ')
Set-TextValue "D37" ('0.0440')
Set-TextValue "E37" ('  +3.59%  ')
Set-TextValue "E38" ('  +0.33%  ')
Set-TextValue "D39" ('3.10')
Set-TextValue "E39" ('  -3.88%  ')
Set-TextValue "D40" ('2.69')
Set-TextValue "E40" ('  -4.43%  ')
Set-TextValue "D41" ('17.48')
Set-TextValue "E41" ('  +1.77%  ')
Set-TextValue "E42" ('  -1.53%  ')
Set-TextValue "E43" ('  +1.01%  ')
Set-TextValue "D44" ('124.13')
Set-TextValue "E44" ('  +2.99%  ')
Set-TextValue "D45" ('22.20')
Set-TextValue "E45" ('  -2.61%  ')
Set-TextValue "D46" ('0.285')
Set-TextValue "E46" ('  +19.88%  ')
Set-TextValue "E47" ('  -3.21%  ')
Set-TextValue "B48" ('ApeXProtocol')
Set-TextValue "C48" ('https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex')
Set-TextValue "D48" ('2.37')
Set-TextValue "E48" ('  +3.17%  ')
Set-TextValue "B49" ('Maker')
Set-TextValue "C49" ('https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr')
Set-TextValue "D49" ('2.044.21')
Set-TextValue "E49" ('  -2.15%  ')
Set-TextValue "B50" ('NEARProtocol')
Set-TextValue "C50" ('https://coinranking.com/coin/DCrsaMv68+nearprotocol-near')
Set-TextValue "D50" ('3.26')
Set-TextValue "E50" ('  +1.04%  ')
Set-TextValue "B51" ('RocketPoolETH')
Set-TextValue "C51" ('https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth')
Set-TextValue "D51" ('3.277.03')
Set-TextValue "E51" ('  +1.48%  ')
